# Parse Excel as Map of Maps #165
# Adds a "Yahoo Home" screen/test entry to the Master and Smoke sheets,
# and the corresponding data row on the Yahoo sheet.

$wb = $excel.ActiveWorkbook

# --- Yahoo sheet: add row 5 (duplicate of row 2, with a new Test name) ---
# Doing this first makes "Yahoo::Yahoo Home" get registered in the shared
# strings table before "Yahoo Home" (matches original authoring order).
$wsYahoo = $wb.Worksheets.Item("Yahoo")
$wsYahoo.Range("A2:G2").Copy($wsYahoo.Range("A5:G5"))
$wsYahoo.Range("A5").Value = "Yahoo::Yahoo Home"
$wsYahoo.Range("B5").Value = ""

# --- Master sheet: add row 4 (duplicate of row 3, pointing at Yahoo Home) ---
$wsMaster = $wb.Worksheets.Item("Master")
$wsMaster.Range("A3:D3").Copy($wsMaster.Range("A4:D4"))
$wsMaster.Range("A4").Value = ""
$wsMaster.Range("B4").Value = "Yahoo"
$wsMaster.Range("C4").Value = "Yahoo Home"
$wsMaster.Range("D4").Value = "Yes"

# --- Smoke sheet: add row 7 (duplicate of row 6, pointing at Yahoo Home) ---
$wsSmoke = $wb.Worksheets.Item("Smoke")
$wsSmoke.Range("A6:D6").Copy($wsSmoke.Range("A7:D7"))
$wsSmoke.Range("A7").Value = ""
$wsSmoke.Range("B7").Value = "Yahoo"
$wsSmoke.Range("C7").Value = "Yahoo Home"
$wsSmoke.Range("D7").Value = "Yes"

# --- Update sheet selections / active-view state ---
$wsSmoke.Range("C11").Select()
$wsYahoo.Range("E11").Select()

# Master becomes the active tab/sheet, selected last so tabSelected sticks.
$wsMaster.Range("C10").Select()
